# Update countries & provincias Spain
# Daily COVID data refresh: updates case counts for several countries and,
# because the sheet is kept sorted by "Casos totales" (column B) descending,
# a handful of neighbouring rows swap positions (country name moves with
# its own row of data) now that the refreshed totals changed their rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Rows whose data was refreshed but keep their country name/position ---

Set-Row 4  653751 5603 56618 563699 13536 846 33434   # Estados Unidos
Set-Row 6  168941 3786 40164 106607 2936 525 22170    # Italia
Set-Row 8  135663 910  77000 54796  4288 63  3867     # Alemania
Set-Row 12 74193  4801 7089  65461  1854 125 1643     # Turquia
Set-Row 16 28923  544  9271  18604  557  38  1048     # Canada
Set-Row 20 14462  112  8986  5083   232  0   393      # Austria
Set-Row 48 3755   141  215   3344   121  7   196      # Republica Dominicana
Set-Row 50 3444   71   526   2849   35   0   69       # Luxemburgo
Set-Row 63 1700   29   703   990    3    0   7        # Barein
Set-Row 76 1128   37   178   918    14   2   32       # Lituania

# --- Rows where updated totals changed the sort rank, so the country name
#     at that row (column A) changes along with the full data row ---

# Rows 29-30: Polonia / Ecuador swap (Ecuador now ranks above Polonia)
$ws.Cells.Item(29, 1).Value = "Ecuador"
Set-Row 29 8225 367 838 6984 168 15 403
$ws.Cells.Item(30, 1).Value = "Polonia"
Set-Row 30 7918 336 774 6830 160 28 314

# Rows 180-181: Islas Virgenes de los Estados Unidos / Fiyi swap
$ws.Cells.Item(180, 1).Value = "Fiyi"
Set-Row 180 17 1 0 17 0 0 0
$ws.Cells.Item(181, 1).Value = "Islas Virgenes de los Estados Unidos"
Set-Row 181 17 0 0 17 0 0 0

# Rows 190-191: Granada / San Cristobal y Nieves swap
$ws.Cells.Item(190, 1).Value = "San Cristobal y Nieves"
Set-Row 190 14 0 0 14 0 0 0
$ws.Cells.Item(191, 1).Value = "Granada"
Set-Row 191 14 0 0 14 2 0 0

# Rows 196-199: Montserrat, Islas Malvinas, Groenlandia, Islas Turcas y Caicos
# reshuffle - Islas Turcas y Caicos jumps to the top of this tied group.
$ws.Cells.Item(196, 1).Value = "Islas Turcas y Caicos"
Set-Row 196 11 1 0 10 0 0 1
$ws.Cells.Item(197, 1).Value = "Montserrat"
Set-Row 197 11 0 1 10 1 0 0
$ws.Cells.Item(198, 1).Value = "Islas Malvinas"
Set-Row 198 11 0 1 10 0 0 0
$ws.Cells.Item(199, 1).Value = "Groenlandia"
Set-Row 199 11 0 11 0 0 0 0
